# Update "想去人数" (want-to-go count, column F) figures on the "展览"
# sheet and the aggregated "全部类型" sheet, reflecting a refreshed
# data export.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet (rows keyed by their row number; column F)
$ws1.Range("F2").Value = 3079
$ws1.Range("F9").Value = 1062
$ws1.Range("F10").Value = 14943
$ws1.Range("F11").Value = 186
$ws1.Range("F12").Value = 145
$ws1.Range("F13").Value = 505
$ws1.Range("F14").Value = 5938
$ws1.Range("F24").Value = 2958
$ws1.Range("F25").Value = 99
$ws1.Range("F26").Value = 10758
$ws1.Range("F29").Value = 123
$ws1.Range("F30").Value = 3758
$ws1.Range("F32").Value = 69

# 全部类型 sheet (same events, offset by one row vs. 展览)
$ws4.Range("F3").Value = 3079
$ws4.Range("F10").Value = 1062
$ws4.Range("F11").Value = 14943
$ws4.Range("F12").Value = 186
$ws4.Range("F13").Value = 145
$ws4.Range("F14").Value = 505
$ws4.Range("F15").Value = 5938
$ws4.Range("F25").Value = 2958
$ws4.Range("F26").Value = 99
$ws4.Range("F28").Value = 10758
$ws4.Range("F31").Value = 123
$ws4.Range("F32").Value = 3758
$ws4.Range("F34").Value = 69
